$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.149.77'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.922.01'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.01'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5058'
$ws.Range('E7').Value = '  -2.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4040'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08277'
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.113'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.25'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  +3.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.423'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '1.906.41'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.339'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.80'
$ws.Range('E17').Value = '  -2.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001100'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06478'
$ws.Range('E19').Value = '  -3.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.60'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9990'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.982'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').Value = '30.188.19'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.29'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.197'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '22.28'
$ws.Range('E26').Value = '  +5.25%  '
$ws.Range('D27').Value = '2.125.44'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.82'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.371'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.96'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.127'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1045'
$ws.Range('E32').Value = '  -2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.003'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.794'
$ws.Range('E34').Value = '  +5.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02454'
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.418'
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06462'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -2.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.837'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6406'
$ws.Range('E41').Value = '  -1.71%  '
$ws.Range('E42').Value = '  -4.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.217'
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9989'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.30'
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.174'
$ws.Range('E46').Value = '  +4.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6004'
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.648'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.74'
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.12'
$ws.Range('E51').Value = '  -0.53%  '
